$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update company id label and refreshed metrics ---
$ws.Range("B2").Value = "'1"
$ws.Range("D2").Value = -0.435
$ws.Range("E2").Value = 0.322
$ws.Range("G2").Value = -6.763285024154589
$ws.Range("H2").Value = -6.763285024154589
$ws.Range("I2").Value = -7.294685990338165
$ws.Range("J2").Value = -7.294685990338165
$ws.Range("K2").Value = 3.92
$ws.Range("L2").Value = 18.93719806763285
$ws.Range("O2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("U2").Value = 1.64
$ws.Range("V2").Value = 0.0623574144486692
$ws.Range("W2").Value = 1.045333333333333
$ws.Range("X2").Value = 0.06683572093055425
$ws.Range("Y2").Value = 0.978497612402779
$ws.Range("Z2").Value = 0.02966040980083107
$ws.Range("AA2").Value = -0.2163633758418112
$ws.Range("AB2").Value = 0.06650621308893809
$ws.Range("AC2").Value = -0.2828695889307493
$ws.Range("AD2").Value = 0.243
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0.243
$ws.Range("AG2").Value = -1.397
$ws.Range("AH2").Value = 0.009154956108955279
$ws.Range("AI2").Value = 0.09669717469160366
$ws.Range("AJ2").Value = -0.05609765891659638
$ws.Range("AK2").Value = -1.600229095074455
$ws.Range("AL2").Value = 0.081
$ws.Range("AM2").Value = 0.081
$ws.Range("AN2").Value = -0.1547770700636943
$ws.Range("AO2").Value = -18.64197530864197
$ws.Range("AP2").Value = 0.889808917197452
$ws.Range("AQ2").Value = -18.64197530864197

# --- Row 3: refreshed metrics (company_name text unchanged) ---
$ws.Range("D3").Value = -0.435
$ws.Range("E3").Value = 0.322
$ws.Range("G3").Value = -6.763285024154589
$ws.Range("H3").Value = -6.763285024154589
$ws.Range("I3").Value = -7.294685990338165
$ws.Range("J3").Value = -7.294685990338165
$ws.Range("K3").Value = 3.92
$ws.Range("L3").Value = 18.93719806763285
$ws.Range("O3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("U3").Value = 1.64
$ws.Range("V3").Value = 0.0623574144486692
$ws.Range("W3").Value = 1.045333333333333
$ws.Range("X3").Value = 0.06683572093055425
$ws.Range("Y3").Value = 0.978497612402779
$ws.Range("Z3").Value = 0.02966040980083107
$ws.Range("AA3").Value = -0.2163633758418112
$ws.Range("AB3").Value = 0.06650621308893809
$ws.Range("AC3").Value = -0.2828695889307493
$ws.Range("AD3").Value = 0.243
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0.243
$ws.Range("AG3").Value = -1.397
$ws.Range("AH3").Value = 0.009154956108955279
$ws.Range("AI3").Value = 0.09669717469160366
$ws.Range("AJ3").Value = -0.05609765891659638
$ws.Range("AK3").Value = -1.600229095074455
$ws.Range("AL3").Value = 0.081
$ws.Range("AM3").Value = 0.081
$ws.Range("AN3").Value = -0.1547770700636943
$ws.Range("AO3").Value = -18.64197530864197
$ws.Range("AP3").Value = 0.889808917197452
$ws.Range("AQ3").Value = -18.64197530864197

# --- Remove row 4 (Natural Health Farm Holdings Inc. record dropped) ---
$ws.Rows(4).EntireRow.Delete()
